$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Requisitos" list occupies rows 25-27 in columns B and C.
# The entry for LOM3246 (Indicação de Conjunto) moves from the first
# position (row 25) to the last position (row 27), shifting the other
# two requirement rows up by one.

$row25 = "LOB1021 -  Física IV  (Requisito)`n"
$row26 = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"
$row27 = "LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)`n"

$ws.Range("B25").Value = $row25
$ws.Range("C25").Value = $row25

$ws.Range("B26").Value = $row26
$ws.Range("C26").Value = $row26

$ws.Range("B27").Value = $row27
$ws.Range("C27").Value = $row27
